$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 258.61
$ws.Range("I15").Value = 258.61
$ws.Range("K15").Value = 775.83
$ws.Range("M15").Value = -606.83

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 243.72728
$ws.Range("I55").Value = 238.71428
$ws.Range("J55").Value = 252.5
$ws.Range("K55").Value = 238.71428
$ws.Range("L55").Value = 252.5
$ws.Range("M55").Value = -24.71428
$ws.Range("N55").Value = -680.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 19238806
$ws.Range("I76").Value = 50012930
$ws.Range("J76").Value = 4979.25
$ws.Range("K76").Value = 50012930
$ws.Range("L76").Value = 4979.25
$ws.Range("M76").Value = -50012615
$ws.Range("N76").Value = -5609.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 19238806
$ws.Range("I79").Value = 50012930
$ws.Range("J79").Value = 4979.25
$ws.Range("K79").Value = 50012930
$ws.Range("L79").Value = 4979.25
$ws.Range("M79").Value = -50011838
$ws.Range("N79").Value = -7163.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1194.1333
$ws.Range("I129").Value = 361.66666
$ws.Range("J129").Value = 1402.25
$ws.Range("K129").Value = 1084.99998
$ws.Range("L129").Value = 4206.75
$ws.Range("M129").Value = 3915.00002
$ws.Range("N129").Value = -14206.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 687.2
$ws.Range("I135").Value = 339
$ws.Range("K135").Value = 3051
$ws.Range("M135").Value = -516

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2805.8
$ws.Range("I137").Value = 2812.4
$ws.Range("J137").Value = 2789.3
$ws.Range("K137").Value = 8437.200000000001
$ws.Range("L137").Value = 8367.900000000001
$ws.Range("M137").Value = -5887.200000000001
$ws.Range("N137").Value = -13467.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 602562.7
$ws.Range("I61").Value = 480169.47
$ws.Range("J61").Value = 837149.75
$ws.Range("K61").Value = 480169.47
$ws.Range("L61").Value = 837149.75
$ws.Range("M61").Value = -479957.47
$ws.Range("N61").Value = -837573.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4745
$ws.Range("I63").Value = 5333.3335
$ws.Range("J63").Value = 2980
$ws.Range("K63").Value = 5333.3335
$ws.Range("L63").Value = 2980
$ws.Range("M63").Value = -4647.3335
$ws.Range("N63").Value = -4352

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 29800
$ws.Range("J64").Value = 29800
$ws.Range("L64").Value = 29800
$ws.Range("N64").Value = -30296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4745
$ws.Range("I66").Value = 5333.3335
$ws.Range("J66").Value = 2980
$ws.Range("K66").Value = 26666.6675
$ws.Range("L66").Value = 14900
$ws.Range("M66").Value = -23234.6675
$ws.Range("N66").Value = -21764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 29800
$ws.Range("J67").Value = 29800
$ws.Range("L67").Value = 29800
$ws.Range("N67").Value = -31516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 602562.7
$ws.Range("I136").Value = 480169.47
$ws.Range("J136").Value = 837149.75
$ws.Range("K136").Value = 1440508.41
$ws.Range("L136").Value = 2511449.25
$ws.Range("M136").Value = -1437958.41
$ws.Range("N136").Value = -2516549.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2275
$ws.Range("I105").Value = 2666.6667
$ws.Range("J105").Value = 2144.4443
$ws.Range("K105").Value = 2666.6667
$ws.Range("L105").Value = 2144.4443
$ws.Range("M105").Value = -919.6667000000002
$ws.Range("N105").Value = -5638.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2066.6667
$ws.Range("I107").Value = 2066.6667
$ws.Range("K107").Value = 2066.6667
$ws.Range("M107").Value = -146.6667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2277.7646
$ws.Range("I31").Value = 1695.225
$ws.Range("J31").Value = 3109.9644
$ws.Range("K31").Value = 1695.225
$ws.Range("L31").Value = 3109.9644
$ws.Range("M31").Value = -1400.225
$ws.Range("N31").Value = -3699.9644

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2277.7646
$ws.Range("I34").Value = 1695.225
$ws.Range("J34").Value = 3109.9644
$ws.Range("K34").Value = 1695.225
$ws.Range("L34").Value = 3109.9644
$ws.Range("M34").Value = -1493.225
$ws.Range("N34").Value = -3513.9644

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1931.8
$ws.Range("I132").Value = 1080.0322
$ws.Range("K132").Value = 3240.0966
$ws.Range("M132").Value = -710.0966000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5677.4443
$ws.Range("I56").Value = 5677.4443
$ws.Range("K56").Value = 5677.4443
$ws.Range("M56").Value = -5147.4443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 28317
$ws.Range("J120").Value = 28317
$ws.Range("L120").Value = 28317
$ws.Range("N120").Value = -37993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3899.8333
$ws.Range("I132").Value = 3768.48
$ws.Range("J132").Value = 4198.364
$ws.Range("K132").Value = 11305.44
$ws.Range("L132").Value = 12595.092
$ws.Range("M132").Value = -8775.440000000001
$ws.Range("N132").Value = -17655.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1497.375
$ws.Range("I61").Value = 1497.375
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1497.375
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1295.375
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1497.375
$ws.Range("I113").Value = 1497.375
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1497.375
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 672.625
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3115.3408
$ws.Range("I132").Value = 2489.95
$ws.Range("J132").Value = 3636.5
$ws.Range("K132").Value = 7469.849999999999
$ws.Range("L132").Value = 10909.5
$ws.Range("M132").Value = -4939.849999999999
$ws.Range("N132").Value = -15969.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 889.5625
$ws.Range("I126").Value = 373.83334
$ws.Range("K126").Value = 1121.50002
$ws.Range("M126").Value = 1348.49998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2731.1155
$ws.Range("I132").Value = 1911.8235
$ws.Range("J132").Value = 4278.6665
$ws.Range("K132").Value = 5735.470499999999
$ws.Range("L132").Value = 12835.9995
$ws.Range("M132").Value = -3205.470499999999
$ws.Range("N132").Value = -17895.9995
